$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Warung ini mempunyai konsep recycle dan reduce yang sangat baik, sangat sesuai dengan program pemerintah tentang pengelolaan sampah yang baik.'
$ws.Range("B2").Value = 'positive'

$ws.Range("A3").Value = 'Sistem pengelolaan sampah ini perlu lebih transparan agar masyarakat lebih percaya diri dalam menggunakannya.'
$ws.Range("B3").Value = 'neutral'

$ws.Range("A4").Value = 'Lokasi tempat sampahnya sangat strategis, membuat pengalaman saya membuang sampah jadi lebih mudah dan efisien.'
$ws.Range("B4").Value = 'positive'

$ws.Range("A5").Value = 'Puas dengan pengalaman penggunaan produk pengelolaan sampah cerdas, membuat hidup jadi lebih mudah dan bersih!'
$ws.Range("B5").Value = 'positive'

$ws.Range("A6").Value = 'Sistem pengelolaan sampah yang buruk membuat kuliah kita jadi tidak nyaman, masih banyak tempat-tempat yang kumuh dan tidak ada inovasi untuk mengatasinya.'
$ws.Range("B6").Value = 'negative'

$ws.Range("A7").Value = 'Pengelolaan sampah yang efektif di food stall ini membuat saya merasa nyaman dan lingkungan tetap bersih, meskipun banyak pengunjung.'
$ws.Range("B7").Value = 'positive'

$ws.Range("A8").Value = 'Saya kecewa dengan sistem pengelolaan sampah cerdas karena tidak membantu saya mengurangi tagihan kartu kredit saya.'
$ws.Range("B8").Value = 'negative'

$ws.Range("A9").Value = 'Lingkungan bersih dan tertata rapi, sangat memuaskan pengalaman makan di sini, apalagi dengan pelayanan prima dan harga sesuai kualitas.'
$ws.Range("B9").Value = 'positive'

$ws.Range("A10").Value = 'Pelayanan pengelolaan limbah di lokasi ini sangat efisien dan memuaskan, membuat kami merasa nyaman dan sangat puas dengan cara mereka mengelola sampah.'
$ws.Range("B10").Value = 'positive'

$ws.Range("A11").Value = 'Merasa nyaman karena sistem pengelolaan sampah cerdas kami membantu menjaga kebersihan lingkungan.'
$ws.Range("B11").Value = 'positive'

$ws.Range("A12").Value = 'Sistem pengelolaan sampah online kita lambat dan tidak responsif, sulit untuk melaporkan insiden sampah bahkan setelah isi pulsa.'
$ws.Range("B12").Value = 'negative'

$ws.Range("A13").Value = 'Sistem pengelolaan sampah yang tidak efektif membuat saya merasa kesal karena tujuan saya untuk memiliki lingkungan yang bersih tidak tercapai.'
$ws.Range("B13").Value = 'negative'

$ws.Range("A14").Value = 'Ketika menggunakan smart waste recycling, saya merasa sangat puas karena sampah saya dapat dikelola dengan baik dan efisien, sehingga membuat lingkungan menjadi lebih bersih dan seimbang.'
$ws.Range("B14").Value = 'positive'

$ws.Range("A15").Value = 'Pengalaman buruk dengan pengumpulan sampah digital hari ini, tidak efektif dan berbelit-belit, membuat saya sangat kecewa dengan layanan tersebut.'
$ws.Range("B15").Value = 'negative'

$ws.Range("A16").Value = 'Kesadaran anak muda terhadap pentingnya pengelolaan sampah cerdas saat ini masih kurang.'
$ws.Range("B16").Value = 'negative'

$ws.Range("A17").Value = 'Sistem kebersihan Jakarta ini sangat tidak efektif, sepertinya Anies dan Sandi tidak tahu cara mengelola sampah yang baik'
$ws.Range("B17").Value = 'negative'

$ws.Range("A18").Value = 'Pembentukan BUMN dengan dasar hukum yang jelas akan meningkatkan kepercayaan masyarakat terhadap pengelolaan sampah yang lebih baik.'
$ws.Range("B18").Value = 'neutral'

$ws.Range("A19").Value = 'Sistem pengelolaan sampah cerdas di kawasan ini sangat membantu mengurangi kekumuhan dan membuat lingkungan menjadi lebih hijau.'
$ws.Range("B19").Value = 'positive'

$ws.Range("A20").Value = 'Pengelolaan sampah elektronik di sini sangat efektif, membuat lingkungan sekitar tetap bersih dan nyaman.'
$ws.Range("B20").Value = 'positive'

$ws.Range("A21").Value = 'Program SGWR 2018 sangat membantu mengajarkan anak-anak untuk menjaga kebersihan lingkungan dan mencintai bumi.'
$ws.Range("B21").Value = 'positive'

